# feat: Add progress #6 'graph-ql'
#
# Adds a new bulleted list item after the paragraph
# "Facebook uses GraphQL since 2012 in their native mobile apps":
#   "GraphQL can be used with any programming language and framework"
# The new paragraph should inherit the same ListParagraph style / numbering
# (ilvl 0, numId 1) as its predecessor.

$d = $word.ActiveDocument

$anchorText = "Facebook uses GraphQL since 2012 in their native mobile apps"
$newText    = "GraphQL can be used with any programming language and framework"

# Locate the paragraph that holds the anchor text.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*$anchorText*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find anchor paragraph containing '$anchorText'"
}

$anchorPara = $d.Paragraphs.Item($targetIndex)

# Insert a new paragraph right after the anchor paragraph. InsertParagraphAfter
# on the paragraph's Range creates a sibling paragraph that inherits the same
# paragraph formatting (pStyle "ListParagraph" + numPr ilvl=0/numId=1).
$anchorPara.Range.InsertParagraphAfter()

# The freshly inserted paragraph is now immediately after the anchor paragraph.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = $newText
